$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("7-17-2022 1000m")

$ws.Range("A3").Value = "Rower 2"
$ws.Range("B3").Value = 220
$ws.Range("C3").Value = 0.0023148148148148151
$ws.Range("D3").Value = 0.0011944444444444446
$ws.Range("E3").Value = 0.0012847222222222223
$ws.Range("F3").Value = 0.0013182870370370371
$ws.Range("G3").Value = 0.001230324074074074
$ws.Range("H3").Value = 0.0010972222222222223
$ws.Range("I3").Value = 0.0013356481481481481

$ws.Range("H4").Select()
